# Refresh computed Leve price/profit columns (H:N) across the crafting-job
# sheets (ALC, BSM, CRP, CUL, GSM, LTW, WVR) with the latest Universalis
# market snapshot, mirroring the scheduled-runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3830.8948
$ws.Range("I64").Value = 3548.9167
$ws.Range("J64").Value = 4314.2856
$ws.Range("K64").Value = 3548.9167
$ws.Range("L64").Value = 4314.2856
$ws.Range("M64").Value = -3300.9167
$ws.Range("N64").Value = -4810.2856
$ws.Range("H67").Value = 3830.8948
$ws.Range("I67").Value = 3548.9167
$ws.Range("J67").Value = 4314.2856
$ws.Range("K67").Value = 3548.9167
$ws.Range("L67").Value = 4314.2856
$ws.Range("M67").Value = -2690.9167
$ws.Range("N67").Value = -6030.2856
$ws.Range("H94").Value = 1953.5
$ws.Range("I94").Value = 1953.5
$ws.Range("K94").Value = 1953.5
$ws.Range("M94").Value = -1502.5
$ws.Range("H113").Value = 3000.1052
$ws.Range("I113").Value = 2096.111
$ws.Range("J113").Value = 3813.7
$ws.Range("K113").Value = 2096.111
$ws.Range("L113").Value = 3813.7
$ws.Range("M113").Value = 1157.889
$ws.Range("N113").Value = -10321.7
$ws.Range("H116").Value = 2183.1667
$ws.Range("I116").Value = 1944.4445
$ws.Range("K116").Value = 1944.4445
$ws.Range("M116").Value = 1497.5555
$ws.Range("H138").Value = 4434.9624
$ws.Range("I138").Value = 1337.0741
$ws.Range("J138").Value = 7652
$ws.Range("K138").Value = 4011.2223
$ws.Range("L138").Value = 22956
$ws.Range("M138").Value = 1128.7777
$ws.Range("N138").Value = -33236

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 841
$ws.Range("I99").Value = 841
$ws.Range("K99").Value = 841
$ws.Range("M99").Value = 657

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 19642.9
$ws.Range("J59").Value = 19642.9
$ws.Range("L59").Value = 19642.9
$ws.Range("N59").Value = -21932.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1666.25
$ws.Range("I22").Value = 260
$ws.Range("J22").Value = 2510
$ws.Range("K22").Value = 780
$ws.Range("L22").Value = 7530
$ws.Range("M22").Value = -611
$ws.Range("N22").Value = -7868
$ws.Range("H27").Value = 1666.25
$ws.Range("I27").Value = 260
$ws.Range("J27").Value = 2510
$ws.Range("K27").Value = 780
$ws.Range("L27").Value = 7530
$ws.Range("M27").Value = -678
$ws.Range("N27").Value = -7734
$ws.Range("H39").Value = 6909.091
$ws.Range("J39").Value = 8975
$ws.Range("L39").Value = 26925
$ws.Range("N39").Value = -27513
$ws.Range("H41").Value = 747.8570999999999
$ws.Range("I41").Value = 530
$ws.Range("J41").Value = 1292.5
$ws.Range("K41").Value = 1590
$ws.Range("L41").Value = 3877.5
$ws.Range("M41").Value = -1252
$ws.Range("N41").Value = -4553.5
$ws.Range("H60").Value = 650.1539
$ws.Range("I60").Value = 273.83334
$ws.Range("J60").Value = 972.7143
$ws.Range("K60").Value = 821.5000200000001
$ws.Range("L60").Value = 2918.1429
$ws.Range("M60").Value = -570.5000200000001
$ws.Range("N60").Value = -3420.1429
$ws.Range("H113").Value = 669
$ws.Range("I113").Value = 697.2432
$ws.Range("J113").Value = 610.94446
$ws.Range("K113").Value = 2091.7296
$ws.Range("L113").Value = 1832.83338
$ws.Range("M113").Value = 78.27039999999988
$ws.Range("N113").Value = -6172.83338
$ws.Range("H131").Value = 1144.0227
$ws.Range("I131").Value = 1210.4286
$ws.Range("J131").Value = 1083.3914
$ws.Range("K131").Value = 3631.2858
$ws.Range("L131").Value = 3250.1742
$ws.Range("M131").Value = 1408.7142
$ws.Range("N131").Value = -13330.1742
$ws.Range("H134").Value = 4768.1665
$ws.Range("I134").Value = 4658.385
$ws.Range("K134").Value = 13975.155
$ws.Range("M134").Value = -8905.155000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 15656.333
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 15656.333
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 15656.333
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -16216.333
$ws.Range("H50").Value = 15656.333
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 15656.333
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 15656.333
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -16652.333
$ws.Range("H58").Value = 19800
$ws.Range("J58").Value = 19800
$ws.Range("L58").Value = 19800
$ws.Range("N58").Value = -20354
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32216
$ws.Range("H126").Value = 2800
$ws.Range("I126").Value = 1955.5555
$ws.Range("J126").Value = 3384.6155
$ws.Range("K126").Value = 5866.666499999999
$ws.Range("L126").Value = 10153.8465
$ws.Range("M126").Value = -3396.666499999999
$ws.Range("N126").Value = -15093.8465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4777.231
$ws.Range("I7").Value = 3819.4285
$ws.Range("J7").Value = 8800
$ws.Range("K7").Value = 3819.4285
$ws.Range("L7").Value = 8800
$ws.Range("M7").Value = -3707.4285
$ws.Range("N7").Value = -9024
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("H40").Value = 3796.25
$ws.Range("J40").Value = 3728.3333
$ws.Range("L40").Value = 3728.3333
$ws.Range("N40").Value = -4000.3333
$ws.Range("H61").Value = 410415.72
$ws.Range("I61").Value = 15148.5625
$ws.Range("J61").Value = 1113112.9
$ws.Range("K61").Value = 15148.5625
$ws.Range("L61").Value = 1113112.9
$ws.Range("M61").Value = -14946.5625
$ws.Range("N61").Value = -1113516.9
$ws.Range("H113").Value = 410415.72
$ws.Range("I113").Value = 15148.5625
$ws.Range("J113").Value = 1113112.9
$ws.Range("K113").Value = 15148.5625
$ws.Range("L113").Value = 1113112.9
$ws.Range("M113").Value = -12978.5625
$ws.Range("N113").Value = -1117452.9
$ws.Range("H122").Value = 6114.0713
$ws.Range("I122").Value = 5525.657
$ws.Range("J122").Value = 7094.7617
$ws.Range("K122").Value = 16576.971
$ws.Range("L122").Value = 21284.2851
$ws.Range("M122").Value = -14126.971
$ws.Range("N122").Value = -26184.2851
$ws.Range("H123").Value = 66900
$ws.Range("J123").Value = 66900
$ws.Range("L123").Value = 66900
$ws.Range("N123").Value = -76700
$ws.Range("H126").Value = 4777.231
$ws.Range("I126").Value = 3819.4285
$ws.Range("J126").Value = 8800
$ws.Range("K126").Value = 11458.2855
$ws.Range("L126").Value = 26400
$ws.Range("M126").Value = -8988.2855
$ws.Range("N126").Value = -31340
$ws.Range("H137").Value = 49997.5
$ws.Range("J137").Value = 49997.5
$ws.Range("L137").Value = 49997.5
$ws.Range("N137").Value = -60197.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 49500
$ws.Range("J115").Value = 49500
$ws.Range("L115").Value = 49500
$ws.Range("N115").Value = -52634
$ws.Range("H126").Value = 1886.2858
$ws.Range("I126").Value = 1886.2858
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5658.857400000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3188.857400000001
$ws.Range("N126").ClearContents()
$ws.Range("H127").Value = 63000
$ws.Range("J127").Value = 63000
$ws.Range("L127").Value = 63000
$ws.Range("N127").Value = -72920
